# Apply the memmap.xlsx "power distribution" rework:
#  - Row 29 (was RTC @ #B4-#B8, IRQ "3 (... 1, 2, or 3)") becomes the first
#    PIT row @ #B4-#B7, IRQ "3 (... /NMI, /INT, 0, or 3)".
#  - Row 30 (was DIAG @ #B8, N/A) becomes a second PIT row @ #B4-#B7,
#    IRQ "3 (... /INT, 5, 6, or 7)", description "Real Time Clock card".
#  - A new row 31 is added for DIAG @ #B8-#BF, N/A, with an updated
#    description noting it does not decode A0-A2.
#  - The lone footer cell in row 35 (col F) moves down to row 36.
#  - Column E is widened to fit the new, longer IRQ description text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29: RTC -> PIT (first range) ---------------------------------
$ws.Range("C29").Value = "#B4-#B7"
$ws.Range("D29").Value = "PIT"
$ws.Range("E29").Value = "3 (but configurable as /NMI, /INT, 0, or 3)"
# F29 ("Real Time Clock card") is unchanged.

# --- Row 30: DIAG -> PIT (second range) --------------------------------
$ws.Range("C30").Value = "#B4-#B7"
$ws.Range("D30").Value = "PIT"
$ws.Range("E30").Value = "3 (but configurable as /INT, 5, 6, or 7)"
$ws.Range("F30").Value = "Real Time Clock card"

# --- Row 31 (new): DIAG, now at #B8-#BF ---------------------------------
$ws.Range("A31").Value = "IO"
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").HorizontalAlignment = -4152
$ws.Range("B31").Value = "*"
$ws.Range("C31").Value = "#B8-#BF"
$ws.Range("D31").Value = "DIAG"
$ws.Range("E31").Value = "N/A"
$ws.Range("F31").Value = "Diagnostics card (does not decode A0-A2)"

# --- Move the lone footer cell from row 35 to row 36 --------------------
$ws.Range("F36").Value = $ws.Range("F35").Value2
$ws.Range("F35").ClearContents()

# --- Update selection to mirror the author's saved cursor position ------
$ws.Range("E31").Select()

# --- Widen column E to fit the new (longer) IRQ description text -------
$ws.Columns.Item(5).ColumnWidth = 34.19029388403495
